$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the formatting of row 14 (C/D/E) to match the rest of the table body ---
$ws.Range("C13:E13").Copy()
$ws.Range("C14:E14").PasteSpecial(-4122)

# --- Set row 14 data (was an empty placeholder row before) ---
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "28"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "35"

# --- Build new rows 15-24 by copying the fully formatted row 14 downward ---
$ws.Range("B14:F14").Copy()
$ws.Range("B15:F24").PasteSpecial(-4122)
$ws.Range("B14:F14").Copy()
for ($r = 15; $r -le 24; $r++) {
  $ws.Range("B" + $r + ":F" + $r).PasteSpecial(-4163)
}

# --- Fill in B (index), C, E, F values for rows 15-24 ---
$data = @(
  @(15, 11, "24", "07", "v3.4"),
  @(16, 12, "21", "03", "v3.4"),
  @(17, 13, "20", "39", "v3.4"),
  @(18, 14, "17", "25", "v3.4"),
  @(19, 15, "18", "23", "v3.4"),
  @(20, 16, "16", "47", "v3.4"),
  @(21, 17, "18", "00", "v3.4"),
  @(22, 18, "17", "39", "v3.4"),
  @(23, 19, "", "", "v3.4"),
  @(24, 20, "", "", "v3.4")
)

foreach ($row in $data) {
  $r = $row[0]
  $idx = $row[1]
  $cval = $row[2]
  $eval = $row[3]
  $fval = $row[4]

  $ws.Range("B$r").Value = $idx
  $ws.Range("C$r").NumberFormat = "@"
  $ws.Range("E$r").NumberFormat = "@"
  if ($cval -ne "") {
    $ws.Range("C$r").Value = $cval
  } else {
    $ws.Range("C$r").ClearContents() | Out-Null
  }
  if ($eval -ne "") {
    $ws.Range("E$r").Value = $eval
  } else {
    $ws.Range("E$r").ClearContents() | Out-Null
  }
  $ws.Range("F$r").Value = $fval
}

# --- Update the view (scrolled position + active cell selection) ---
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("H22").Select() | Out-Null
